# Sensitivity analysis: mark "Crisis" (column G) = 1 for the selected
# deal/firm-characteristic observations that occurred during a crisis period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(31, 37, 38, 48, 49, 50, 51, 52, 54, 55, 56, 57, 59, 60, 61, 62, 76, 77, `
          113, 116, 117, 118, 119, 120, 121, 122, 123, 124, 174, 180, 181, 182, `
          183, 184, 185, 186, 187, 222, 223, 224, 225, 226, 227, 228, 229, 230, `
          231, 232, 233, 234, 235, 236, 237, 238, 239, 240, 241, 242, 243, 244, 245)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 1
}
